$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '76.391.07'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.076.61'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.30%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '198.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '619.94'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.45%  '
$ws.Range('E8').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E9').ClearFormats()

$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Cardano'
$ws.Range('B10').ClearFormats()
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('C10').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.449'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('E10').ClearFormats()

$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'TRON'
$ws.Range('B11').ClearFormats()
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C11').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.161'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E11').ClearFormats()

$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('B12').ClearFormats()
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C12').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.25'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.07%  '
$ws.Range('E12').ClearFormats()

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B13').ClearFormats()
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C13').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.641.84'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('E13').ClearFormats()

$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('B14').ClearFormats()
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C14').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.29'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.97%  '
$ws.Range('E14').ClearFormats()

$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('B15').ClearFormats()
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C15').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000200'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.37%  '
$ws.Range('E15').ClearFormats()

$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('B16').ClearFormats()
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C16').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.286.07'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('E16').ClearFormats()

$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('B17').ClearFormats()
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C17').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.069.89'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.23%  '
$ws.Range('E17').ClearFormats()

$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('B18').ClearFormats()
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C18').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.47'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('E18').ClearFormats()

$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('B19').ClearFormats()
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C19').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.02'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.16%  '
$ws.Range('E19').ClearFormats()

$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('B20').ClearFormats()
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('C20').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +16.02%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '387.38'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('E21').ClearFormats()

$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('B22').ClearFormats()
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C22').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.52'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.30%  '
$ws.Range('E22').ClearFormats()

$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'LEO'
$ws.Range('B23').ClearFormats()
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C23').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.47'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('E23').ClearFormats()

$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'NEARProtocol'
$ws.Range('B24').ClearFormats()
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C24').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.56'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.00%  '
$ws.Range('E24').ClearFormats()

$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('B25').ClearFormats()
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('C25').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.231.39'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.44%  '
$ws.Range('E25').ClearFormats()

$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('B26').ClearFormats()
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C26').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('E26').ClearFormats()

$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Dai'
$ws.Range('B27').ClearFormats()
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C27').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E27').ClearFormats()

$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Aptos'
$ws.Range('B28').ClearFormats()
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C28').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.40%  '
$ws.Range('E28').ClearFormats()

$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'PEPE'
$ws.Range('B29').ClearFormats()
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C29').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000109'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('E29').ClearFormats()

$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('B30').ClearFormats()
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('C30').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E30').ClearFormats()

$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('B31').ClearFormats()
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C31').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.32'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('E31').ClearFormats()

$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('B32').ClearFormats()
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C32').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.43'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('E32').ClearFormats()

$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('B33').ClearFormats()
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C33').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '501.62'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('E33').ClearFormats()

$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('B34').ClearFormats()
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C34').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.93'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.80%  '
$ws.Range('E34').ClearFormats()

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('B35').ClearFormats()
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C35').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +18.68%  '
$ws.Range('E35').ClearFormats()

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('B36').ClearFormats()
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C36').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E36').ClearFormats()

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('B37').ClearFormats()
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C37').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.88'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.09%  '
$ws.Range('E37').ClearFormats()

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Monero'
$ws.Range('B38').ClearFormats()
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C38').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.56'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('E38').ClearFormats()

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Aave'
$ws.Range('B39').ClearFormats()
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C39').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '194.16'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.45%  '
$ws.Range('E39').ClearFormats()

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('B40').ClearFormats()
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('C40').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '20.07'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E40').ClearFormats()

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.378'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('E41').ClearFormats()

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Cronos'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.102'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.07%  '
$ws.Range('E42').ClearFormats()

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'USDe'
$ws.Range('B43').ClearFormats()
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C43').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('E43').ClearFormats()

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Mantle'
$ws.Range('B44').ClearFormats()
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C44').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.799'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +20.26%  '
$ws.Range('E44').ClearFormats()

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.19'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.03%  '
$ws.Range('E45').ClearFormats()

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('B46').ClearFormats()
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C46').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.26'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.16%  '
$ws.Range('E46').ClearFormats()

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Stacks'
$ws.Range('B47').ClearFormats()
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C47').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.66'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('E47').ClearFormats()

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('B48').ClearFormats()
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C48').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.47'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.49%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '40.99'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('E49').ClearFormats()

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('B50').ClearFormats()
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C50').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.600'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('E50').ClearFormats()

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('B51').ClearFormats()
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C51').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.93'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.31%  '
$ws.Range('E51').ClearFormats()
